# corrige calculo semanas de subgrupos
# Recalcula la semana/subgrupo (A/B/C) asignada a varios alumnos en las
# columnas subgrupo_potencia (C) y subgrupo_robotica (D) para corregir
# colisiones de sesion/semana.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "-"
$ws.Range("D5").Value = "-"
$ws.Range("C9").Value = "MI11-A"
$ws.Range("C17").Value = "MI11-B"
$ws.Range("C18").Value = "MI11-B"
$ws.Range("D22").Value = "-"
$ws.Range("C23").Value = "MI11-B"
$ws.Range("D35").Value = "MA11-A"
$ws.Range("D55").Value = "-"
$ws.Range("D62").Value = "-"
$ws.Range("C63").Value = "MI11-C"
$ws.Range("C75").Value = "MI11-C"
$ws.Range("D75").Value = "MI11-B"
$ws.Range("D79").Value = "MA11-A"
$ws.Range("C90").Value = "JU11-A"
$ws.Range("D119").Value = "MA11-B"
$ws.Range("C120").Value = "JU11-B"
$ws.Range("C135").Value = "JU11-C"
$ws.Range("D146").Value = "MA11-B"
$ws.Range("D149").Value = "MA11-C"
$ws.Range("D150").Value = "MA11-B"
$ws.Range("D163").Value = "MA11-C"
$ws.Range("C197").Value = "JU11-D"
$ws.Range("D231").Value = "MA11-A"
$ws.Range("D242").Value = "MA11-C"
